$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45896, 5341.91530752175, 4374.23088654824, 1944, 5689.101763, 115.725722584437),
    @(45897, 5341.91604070457, 4350.85613728879, 1944, 5689.102522, 114.751775774342),
    @(45898, 5341.91604070659, 3708.10312613918, 1944, 5689.102522, 87.9704003096911),
    @(45899, 975.205470885634, 1610.97852240698, 1944, 1881.071313, 23.8685151883894),
    @(45900, 862.92289373491, 1568.06896867702, 1944, 1758.699924, 21.6602499559213),
    @(45901, 5776.55841866516, 4973.0316110909, 2952, 6391.159748, 109.818039184406),
    @(45902, 5784.11073458408, 5172.02032529119, 2952, 6398.77612, 118.111904612796),
    @(45903, 5784.11073458408, 5204.06747854613, 2952, 6398.77612, 119.447202665086),
    @(45904, 5812.79253171584, 5210.12562566099, 2952, 6445.468018, 120.450046331048),
    @(45905, 5812.79253171584, 4530.42391900527, 2952, 6445.468018, 92.1291418870597),
    @(45906, 1096.28225605742, 2275.43932428577, 2952, 2113.128692, 14.1785733428476),
    @(45907, 970.493609902267, 2230.88381429177, 2952, 1974.388171, 11.7824323078961),
    @(45908, 6110.60514061168, 5077.61908918156, 2952, 6751.185399, 115.258306148745),
    @(45909, 6110.60514061168, 5260.86093663226, 2952, 6751.185399, 122.893383125858)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
